# Update the passive-trial sample values (columns B:E, rows 1-3) on sheet "Ark1"
# with the new re-exported data, matching the "Hjemme passive tweaks lichtwark
# deleted values" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ark1")

# Row 1 - sample size headers
$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

# Row 2 - "CON" data
$ws.Range("B2").Value = 48.25441703816243
$ws.Range("C2").Value = 56.772752221374418
$ws.Range("D2").Value = 51.273337329918661
$ws.Range("E2").Value = 56.865557802866988

# Row 3 - "STR" data
$ws.Range("B3").Value = 45.78411335805194
$ws.Range("C3").Value = 46.858579123615733
$ws.Range("D3").Value = 45.914231469102674
$ws.Range("E3").Value = 55.734007222601392

# The author's selection now only spans the edited block instead of the
# whole table (B1:AY3 -> B1:E3).
$ws.Range("B1:E3").Select()
